$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the event title (cell B3, merged B3:J3) to reflect the new event.
$ws.Range("B3").Value = "Teste de Evento 2 - 2022-09-26 10:15:15"

# Remove the second registrant row (row 6: cadu / cadedu@gmail.com / ...).
$ws.Rows("6:6").Delete() | Out-Null

# Keep the selection on the new last data row, matching the sheet's prior convention.
$ws.Range("B5:J5").Select() | Out-Null
